$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.224.38'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '1.904.23'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5409'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3808'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07300'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9036'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08176'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.42'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.345'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.9987'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008655'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '1.316.22'
$ws.Range("E18").Value = '  -30.76%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '27.247.23'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.046'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.515'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.304'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.751'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.851'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.664'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09197'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8248'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05070'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.224'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.015'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.317'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.708'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6033'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02001'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.58%  '
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.270'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.674'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5184'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '115.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1532'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.85%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9985'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.641'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06090'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
